$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "Transaction #" column (G) for the rows that were missing it
$ws.Range("G2").Value = 732837
$ws.Range("G3").Value = 717203
$ws.Range("G6").Value = 717203
$ws.Range("G7").Value = 732838
$ws.Range("G8").Value = 717203

# Add a comment on row 5 (Comment column J)
$ws.Range("J5").Value = "Transaction has been processed"
